# Add 9 new vocabulary entries (FOREIGN / ENGLISH / DATE) to the "words"
# sheet, appended after the existing 49 data rows (new rows 50-58), all
# dated 2020-12-02 - mirrors an "updated DBs" batch import.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newWords = @(
    "即",
    "热量",
    "配比",
    "摆放",
    "摄入",
    "对应",
    "除此以外",
    "逐渐",
    "养成"
)

$newDefs = @(
    "1.) v. 就是; 2.) adv. Immediately",
    "1.) n. calories; 2.) heat, amount of heat",
    "1.) matching, proportioning",
    "1.) v. to set up, arrange, lay out",
    "1.) v. to take in, to absorb, to consume; 2.) intake, consumption",
    "1.) v. to correspond; 2.) adj. corresponding",
    "1.) prep. beside, except for; 2.) prep. phr. in addition to; 3.) adv. otherwise",
    "1.) adv. gradually, little by little",
    "1.) v. to cultivate, to raise, to form [a habit], to acquire"
)

$newDate = "2020-12-02"

$startRow = 50
$endRow = $startRow + $newWords.Length - 1

# Pre-format the DATE column as text so Excel doesn't silently convert the
# "yyyy-mm-dd" strings into date serials (the source data stores dates as
# plain text, matching every other row in the sheet).
$dateRange = $ws.Range("C$startRow`:C$endRow")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $newWords.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newWords[$i]
    $ws.Cells.Item($r, 2).Value = $newDefs[$i]
    $ws.Cells.Item($r, 3).Value = $newDate
}

# Strip the temporary text-format style back to the sheet's default so the
# new cells stay visually identical to the rest of the table.
$dateRange.Style = "Normal"
